$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 309.375
$ws.Range("I28").Value = 234.5238
$ws.Range("J28").Value = 833.3333
$ws.Range("K28").Value = 234.5238
$ws.Range("L28").Value = 833.3333
$ws.Range("M28").Value = 250.4762
$ws.Range("N28").Value = -1803.3333
$ws.Range("H121").Value = 1280.5
$ws.Range("I121").Value = 561
$ws.Range("J121").Value = 2000
$ws.Range("K121").Value = 1683
$ws.Range("L121").Value = 6000
$ws.Range("M121").Value = 64
$ws.Range("N121").Value = -9494
$ws.Range("H135").Value = 567.7
$ws.Range("I135").Value = 576.931
$ws.Range("J135").Value = 300
$ws.Range("K135").Value = 5192.379000000001
$ws.Range("L135").Value = 2700
$ws.Range("M135").Value = -2657.379000000001
$ws.Range("N135").Value = -7770
$ws.Range("H141").Value = 778799.3
$ws.Range("I141").Value = 1687.1482
$ws.Range("J141").Value = 4275804
$ws.Range("K141").Value = 5061.444600000001
$ws.Range("L141").Value = 12827412
$ws.Range("M141").Value = 118.5553999999993
$ws.Range("N141").Value = -12837772

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2972.71
$ws.Range("I32").Value = 2338.1428
$ws.Range("J32").Value = 9388.888999999999
$ws.Range("K32").Value = 2338.1428
$ws.Range("L32").Value = 9388.888999999999
$ws.Range("M32").Value = -2051.1428
$ws.Range("N32").Value = -9962.888999999999
$ws.Range("H43").Value = 20000
$ws.Range("I43").Value = 20000
$ws.Range("K43").Value = 20000
$ws.Range("M43").Value = -19687
$ws.Range("H61").Value = 2425.9333
$ws.Range("I61").Value = 1615.75
$ws.Range("K61").Value = 1615.75
$ws.Range("M61").Value = -1403.75
$ws.Range("H132").Value = 20836932
$ws.Range("I132").Value = 26318926
$ws.Range("J132").Value = 5352.6
$ws.Range("K132").Value = 78956778
$ws.Range("L132").Value = 16057.8
$ws.Range("M132").Value = -78954248
$ws.Range("N132").Value = -21117.8
$ws.Range("H136").Value = 2425.9333
$ws.Range("I136").Value = 1615.75
$ws.Range("K136").Value = 4847.25
$ws.Range("M136").Value = -2297.25

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3042.375
$ws.Range("I20").Value = 1337.5
$ws.Range("J20").Value = 4747.25
$ws.Range("K20").Value = 1337.5
$ws.Range("L20").Value = 4747.25
$ws.Range("M20").Value = -1090.5
$ws.Range("N20").Value = -5241.25
$ws.Range("H105").Value = 1497.0714
$ws.Range("I105").Value = 1278.421
$ws.Range("J105").Value = 1958.6666
$ws.Range("K105").Value = 1278.421
$ws.Range("L105").Value = 1958.6666
$ws.Range("M105").Value = 468.579
$ws.Range("N105").Value = -5452.6666
$ws.Range("H134").Value = 8784.416999999999
$ws.Range("I134").Value = 9342.714
$ws.Range("J134").Value = 8002.8
$ws.Range("K134").Value = 28028.142
$ws.Range("L134").Value = 24008.4
$ws.Range("M134").Value = -25493.142
$ws.Range("N134").Value = -29078.4

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2328860
$ws.Range("I31").Value = 3127093.2
$ws.Range("J31").Value = 6727.1816
$ws.Range("K31").Value = 3127093.2
$ws.Range("L31").Value = 6727.1816
$ws.Range("M31").Value = -3126798.2
$ws.Range("N31").Value = -7317.1816
$ws.Range("H34").Value = 2328860
$ws.Range("I34").Value = 3127093.2
$ws.Range("J34").Value = 6727.1816
$ws.Range("K34").Value = 3127093.2
$ws.Range("L34").Value = 6727.1816
$ws.Range("M34").Value = -3126891.2
$ws.Range("N34").Value = -7131.1816
$ws.Range("H58").Value = 33337218
$ws.Range("I58").Value = 1575.8
$ws.Range("J58").Value = 50005040
$ws.Range("K58").Value = 1575.8
$ws.Range("L58").Value = 50005040
$ws.Range("M58").Value = -1372.8
$ws.Range("N58").Value = -50005446
$ws.Range("H99").Value = 2648.111
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 3222.1667
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 3222.1667
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -6218.1667
$ws.Range("H106").Value = 23828.4
$ws.Range("J106").Value = 23828.4
$ws.Range("L106").Value = 23828.4
$ws.Range("N106").Value = -26352.4
$ws.Range("H126").Value = 2648.111
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 3222.1667
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 9666.500100000001
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -14606.5001
$ws.Range("H134").Value = 1365.5312
$ws.Range("I134").Value = 819
$ws.Range("J134").Value = 2408.9092
$ws.Range("K134").Value = 2457
$ws.Range("L134").Value = 7226.7276
$ws.Range("M134").Value = 78
$ws.Range("N134").Value = -12296.7276
$ws.Range("H136").Value = 33337218
$ws.Range("I136").Value = 1575.8
$ws.Range("J136").Value = 50005040
$ws.Range("K136").Value = 4727.4
$ws.Range("L136").Value = 150015120
$ws.Range("M136").Value = -2177.4
$ws.Range("N136").Value = -150020220

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1815.3572
$ws.Range("I5").Value = 571.5
$ws.Range("K5").Value = 1714.5
$ws.Range("M5").Value = -1602.5
$ws.Range("H102").Value = 2500
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 3000
$ws.Range("L102").Value = 9000
$ws.Range("M102").Value = -566
$ws.Range("N102").Value = -13868
$ws.Range("H110").Value = 2776
$ws.Range("I110").Value = 2475
$ws.Range("J110").Value = 3980
$ws.Range("K110").Value = 7425
$ws.Range("L110").Value = 11940
$ws.Range("M110").Value = -3335
$ws.Range("N110").Value = -20120
$ws.Range("H114").Value = 1143.409
$ws.Range("I114").Value = 615.875
$ws.Range("J114").Value = 1444.8572
$ws.Range("K114").Value = 1847.625
$ws.Range("L114").Value = 4334.571599999999
$ws.Range("M114").Value = 1406.375
$ws.Range("N114").Value = -10842.5716
$ws.Range("H124").Value = 34400
$ws.Range("I124").Value = 1600
$ws.Range("J124").Value = 100000
$ws.Range("K124").Value = 4800
$ws.Range("L124").Value = 300000
$ws.Range("M124").Value = 110
$ws.Range("N124").Value = -309820
$ws.Range("H131").Value = 1419.8918
$ws.Range("J131").Value = 1293.1852
$ws.Range("L131").Value = 3879.5556
$ws.Range("N131").Value = -13959.5556
$ws.Range("H132").Value = 4855.8
$ws.Range("J132").Value = 6025
$ws.Range("L132").Value = 54225
$ws.Range("N132").Value = -59285
$ws.Range("H135").Value = 1815.3572
$ws.Range("I135").Value = 571.5
$ws.Range("K135").Value = 5143.5
$ws.Range("M135").Value = -2608.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4186.4287
$ws.Range("I70").Value = 4008.75
$ws.Range("J70").Value = 4423.3335
$ws.Range("K70").Value = 4008.75
$ws.Range("L70").Value = 4423.3335
$ws.Range("M70").Value = -3738.75
$ws.Range("N70").Value = -4963.3335
$ws.Range("H73").Value = 4186.4287
$ws.Range("I73").Value = 4008.75
$ws.Range("J73").Value = 4423.3335
$ws.Range("K73").Value = 4008.75
$ws.Range("L73").Value = 4423.3335
$ws.Range("M73").Value = -3072.75
$ws.Range("N73").Value = -6295.3335
$ws.Range("H132").Value = 3165.9143
$ws.Range("J132").Value = 4153.6924
$ws.Range("L132").Value = 12461.0772
$ws.Range("N132").Value = -17521.0772
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = 0

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 62501908
$ws.Range("I61").Value = 71429896
$ws.Range("K61").Value = 71429896
$ws.Range("M61").Value = -71429694
$ws.Range("H68").Value = 1529.15
$ws.Range("I68").Value = 1005.3333
$ws.Range("J68").Value = 3100.6
$ws.Range("K68").Value = 1005.3333
$ws.Range("L68").Value = 3100.6
$ws.Range("M68").Value = -256.3333
$ws.Range("N68").Value = -4598.6
$ws.Range("H71").Value = 1529.15
$ws.Range("I71").Value = 1005.3333
$ws.Range("J71").Value = 3100.6
$ws.Range("K71").Value = 5026.6665
$ws.Range("L71").Value = 15503
$ws.Range("M71").Value = -1282.6665
$ws.Range("N71").Value = -22991
$ws.Range("H113").Value = 62501908
$ws.Range("I113").Value = 71429896
$ws.Range("K113").Value = 71429896
$ws.Range("M113").Value = -71427726
$ws.Range("H122").Value = 2741.9807
$ws.Range("I122").Value = 2434.4285
$ws.Range("J122").Value = 4033.7
$ws.Range("K122").Value = 7303.2855
$ws.Range("L122").Value = 12101.1
$ws.Range("M122").Value = -4853.2855
$ws.Range("N122").Value = -17001.1
$ws.Range("H132").Value = 2825.139
$ws.Range("I132").Value = 1744.4375
$ws.Range("J132").Value = 3689.7
$ws.Range("K132").Value = 5233.3125
$ws.Range("L132").Value = 11069.1
$ws.Range("M132").Value = -2703.3125
$ws.Range("N132").Value = -16129.1

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1452846.5
$ws.Range("I5").Value = 250528.5
$ws.Range("J5").Value = 2414700.8
$ws.Range("K5").Value = 250528.5
$ws.Range("L5").Value = 2414700.8
$ws.Range("M5").Value = -250416.5
$ws.Range("N5").Value = -2414924.8
$ws.Range("H45").Value = 15339.8
$ws.Range("J45").Value = 15339.8
$ws.Range("L45").Value = 15339.8
$ws.Range("N45").Value = -16321.8
$ws.Range("H132").Value = 291634.06
$ws.Range("I132").Value = 402655.75
$ws.Range("J132").Value = 14079.8
$ws.Range("K132").Value = 1207967.25
$ws.Range("L132").Value = 42239.39999999999
$ws.Range("M132").Value = -1205437.25
$ws.Range("N132").Value = -47299.39999999999
$ws.Range("H136").Value = 1887.875
$ws.Range("I136").Value = 840.6
$ws.Range("J136").Value = 3633.3333
$ws.Range("K136").Value = 2521.8
$ws.Range("L136").Value = 10899.9999
$ws.Range("M136").Value = 28.19999999999982
$ws.Range("N136").Value = -15999.9999
